# Append a freshly scraped job listing (2025-08-28 01:16 JST) to the top of
# the "ランサーズ" sheet, and log a matching stats row on the "統計" sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "ランサーズ" (job list) ---------------------------------
$ws1 = $wb.Worksheets.Item("ランサーズ")

# Push existing data rows down by inserting a new blank row at row 2.
$ws1.Rows.Item(2).Insert()

# Fill in the new row 2 with the newly scraped listing.
$ws1.Range("A2").Value = "2025-08-28 01:16:20"
$ws1.Range("B2").Value = "【急募】GoogleMAP機能追加・編集の依頼"
$ws1.Range("C2").Value = "システム開発"
$ws1.Range("D2").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws1.Range("E2").Value = "期限情報なし"
$ws1.Range("F2").Style = "Hyperlink"
$ws1.Range("F2").Value = "https://www.lancers.jp/work/detail/5381118"
$ws1.Range("G2").Value = 18

# The row that used to be row 13 (last data row) is now row 14 after the
# insert; it already carries its text/style, but the Insert() operation
# does not recreate its hyperlink relationship, so add it explicitly.
$ws1.Hyperlinks.Add($ws1.Range("F14"), "https://www.lancers.jp/work/detail/5380420")

# --- Sheet 2: "統計" (stats log) --------------------------------------
$ws2 = $wb.Worksheets.Item("統計")

$ws2.Range("A7").Value = "2025-08-28T01:16:20.703321"
$ws2.Range("B7").Value = 13
$ws2.Range("C7").Value = "全案件リスト"
$ws2.Range("D7").Value = 69.2
$ws2.Range("E7").Value = 4
$ws2.Range("F7").Value = 6
$ws2.Range("G7").Value = 13
